$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the specs values on row 2 to the latest Specs and Settlements figures
$ws.Range("H2").Value = 6        # NumPeoplePerHHRural
$ws.Range("I2").Value = 4        # NumPeoplePerHHUrban
$ws.Range("L2").Value = 0.1322   # GridPrice
$ws.Range("R2").Value = 0.29     # ElecActual
$ws.Range("S2").Value = 0.5523046117685825  # ElecModelled

# Reset the view back to the default top-left / selection at A1
$ws.Range("A1").Select()

$wb.Save()
